# Ch 13 03 - Implementing common DAX query patterns
# Fills in the "Function" reference table (columns C..I) for rows 81-92
# on Sheet1, completing section 13-03 "Implementing common DAX query
# patterns" and appending twelve new G/H/I (Key/Subsection/Function)
# entries. Cells are written in the same order the strings were first
# entered so new shared-string table entries land at matching indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 81: finish 130301 ----------------------------------------------
$ws.Range("H81").Value = "Using ROW to test measures"
$ws.Range("I81").Value = "ROW, CALCULATETABLE"

# --- Row 82: 130302 -------------------------------------------------------
$ws.Range("E82").Value = 3
$ws.Range("F82").Value = "Implementing common DAX query patterns"
$ws.Range("G82").Value = 2
$ws.Range("I82").Value = "SUMMARIZE"

# --- Row 83: 130303 -------------------------------------------------------
$ws.Range("E83").Value = 3
$ws.Range("F83").Value = "Implementing common DAX query patterns"
$ws.Range("G83").Value = 3
$ws.Range("I83").Value = "SUMMARIZECOLUMNS, ROLLUPADDISSUBTOTAL, ROLLUPGROUP, FILTER"

# --- Row 84: 130304 -------------------------------------------------------
$ws.Range("E84").Value = 3
$ws.Range("F84").Value = "Implementing common DAX query patterns"
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = "Using TOPN"

$ws.Range("H83").Value = "Using SUMMARIZECOLUMNS"
$ws.Range("I84").Value = "TOPN"

# --- Row 82: Subsection (reuses the existing "Using SUMMARIZE" string) ---
$ws.Range("H82").Value = "Using SUMMARIZE"

# --- Row 85: 130005 (Key column only, no Chapter/Section) -----------------
$ws.Range("G85").Value = 5
$ws.Range("H85").Value = "Using GENERATE and GENERATEALL"

# --- Row 86: 130006 -------------------------------------------------------
$ws.Range("G86").Value = 6
$ws.Range("H86").Value = "Using ISONORAFTER"
$ws.Range("I86").Value = "ISONORAFTER"

$ws.Range("I85").Value = "GENERATE, GENERATEALL"

# --- Row 87: 130007 -------------------------------------------------------
$ws.Range("G87").Value = 7
$ws.Range("H87").Value = "Using ADDMISSINGITEMS"
$ws.Range("I87").Value = "ADDMISSINGITEMS"

# --- Row 88: 130008 -------------------------------------------------------
$ws.Range("G88").Value = 8
$ws.Range("H88").Value = "Using TOPNSKIP"
$ws.Range("I88").Value = "TOPNSKIP"

# --- Row 89: Key-only continuation row ------------------------------------
$ws.Range("G89").Value = 9
$ws.Range("G89").NumberFormat = "00"
$ws.Range("I89").Value = "GROUPBY"
$ws.Range("H89").Value = "Using GROUBY"

# --- Row 90: brand-new row -------------------------------------------------
$ws.Range("G90").Value = 10
$ws.Range("G90").NumberFormat = "00"
$ws.Range("I90").Value = "NATURALINNERJOIN, NATURALLEFTOUTERJOIN"

# --- Row 91: brand-new row -------------------------------------------------
$ws.Range("G91").Value = 11
$ws.Range("G91").NumberFormat = "00"
$ws.Range("H91").Value = "Using SUBSTITUTEWITHINDEX"

$ws.Range("H90").Value = "Using NATURALINNERJOIN and NATURALLEFTOUTERJOIN"
$ws.Range("I91").Value = "SUBSTITUTEWITHINDEX"

# --- Row 92: brand-new row -------------------------------------------------
$ws.Range("G92").Value = 12
$ws.Range("G92").NumberFormat = "00"
$ws.Range("H92").Value = "Using SAMPLE"
$ws.Range("I92").Value = "SAMPLE"

# Move the active selection to where the author left off editing
$ws.Range("I93").Select()
